$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update monthly target values on row 2
$ws.Range("F2").Value = 11750000
$ws.Range("I2").Value = 9750000
$ws.Range("L2").Value = 9650000

# Update the active cell selection to I9
$ws.Range("I9").Select()
